$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 3206
$ws.Range("I3").Value = 3300
$ws.Range("B4").Value = 1654
$ws.Range("D4").Value = 1924
$ws.Range("H4").Value = 1669
$ws.Range("I4").Value = 775
$ws.Range("I5").Value = 303
$ws.Range("I6").Value = 3748
$ws.Range("B7").Value = 23286
$ws.Range("D7").Value = 28114
$ws.Range("H7").Value = 25977
$ws.Range("I7").Value = 11332

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I5").Value = 12
$ws.Range("I6").Value = 100
$ws.Range("I7").Value = 365

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I2").Value = 59
$ws.Range("I3").Value = 75
$ws.Range("I7").Value = 212

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 159
$ws.Range("I7").Value = 450

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I2").Value = 37
$ws.Range("I7").Value = 98

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 83
$ws.Range("I3").Value = 66
$ws.Range("I6").Value = 83
$ws.Range("I7").Value = 255

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 103
$ws.Range("I4").Value = 45
$ws.Range("I5").Value = 38
$ws.Range("I6").Value = 78
$ws.Range("I7").Value = 378
$ws.Range("I8").Value = 711
$ws.Range("I10").Value = 82
$ws.Range("I11").Value = 178
$ws.Range("I18").Value = 79
$ws.Range("I19").Value = 303
$ws.Range("I20").Value = 290
$ws.Range("I22").Value = 30
$ws.Range("I23").Value = 102
$ws.Range("I25").Value = 52
$ws.Range("I27").Value = 101
$ws.Range("I29").Value = 737
$ws.Range("I33").Value = 512
$ws.Range("I37").Value = 365
$ws.Range("I40").Value = 21
$ws.Range("I41").Value = 51
$ws.Range("D42").Value = 1217
$ws.Range("I42").Value = 398
$ws.Range("I44").Value = 84
$ws.Range("I45").Value = 22
$ws.Range("I48").Value = 134
$ws.Range("I51").Value = 103
$ws.Range("I52").Value = 241
$ws.Range("I53").Value = 124
$ws.Range("I54").Value = 252
$ws.Range("I55").Value = 126
$ws.Range("I57").Value = 40
$ws.Range("I60").Value = 56
$ws.Range("B63").Value = 360
$ws.Range("H63").Value = 203
$ws.Range("I64").Value = 103
$ws.Range("I65").Value = 255
$ws.Range("I67").Value = 450
$ws.Range("I68").Value = 35
$ws.Range("I71").Value = 31
$ws.Range("I72").Value = 41
$ws.Range("I75").Value = 38
$ws.Range("I76").Value = 177
$ws.Range("I78").Value = 155
$ws.Range("I79").Value = 287
$ws.Range("I81").Value = 9
$ws.Range("I83").Value = 232
$ws.Range("I84").Value = 98
$ws.Range("I90").Value = 141
$ws.Range("I94").Value = 99
$ws.Range("I99").Value = 212
$ws.Range("B101").Value = 23286
$ws.Range("D101").Value = 28114
$ws.Range("H101").Value = 25977
$ws.Range("I101").Value = 11332

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 83
$ws.Range("I7").Value = 232

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 124
$ws.Range("I3").Value = 179
$ws.Range("I6").Value = 167
$ws.Range("I7").Value = 512

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 131
$ws.Range("I7").Value = 252

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 223
$ws.Range("I3").Value = 259
$ws.Range("I5").Value = 29
$ws.Range("I7").Value = 737

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I3").Value = 84
$ws.Range("I6").Value = 84
$ws.Range("I7").Value = 303

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I3").Value = 20
$ws.Range("I6").Value = 27
$ws.Range("I7").Value = 84

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I4").Value = 13
$ws.Range("I6").Value = 75
$ws.Range("I7").Value = 134

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I6").Value = 75
$ws.Range("I7").Value = 177

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I2").Value = 34
$ws.Range("I7").Value = 78

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("I6").Value = 12
$ws.Range("I7").Value = 51

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 137
$ws.Range("D4").Value = 73
$ws.Range("D7").Value = 1217
$ws.Range("I7").Value = 398

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I2").Value = 27
$ws.Range("I7").Value = 82

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I6").Value = 62
$ws.Range("I7").Value = 155

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I3").Value = 33
$ws.Range("I7").Value = 126

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I3").Value = 34
$ws.Range("I7").Value = 102

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 83
$ws.Range("I7").Value = 287

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I3").Value = 34
$ws.Range("I7").Value = 103

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I2").Value = 82
$ws.Range("I3").Value = 85
$ws.Range("I7").Value = 290

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I2").Value = 21
$ws.Range("I7").Value = 79

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I5").Value = 9
$ws.Range("I6").Value = 56
$ws.Range("I7").Value = 241

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I3").Value = 16
$ws.Range("I7").Value = 99

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("I2").Value = 14
$ws.Range("I7").Value = 52

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I6").Value = 41
$ws.Range("I7").Value = 178

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I2").Value = 35
$ws.Range("I7").Value = 103

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I3").Value = 195
$ws.Range("I6").Value = 229
$ws.Range("I7").Value = 711

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("I6").Value = 19
$ws.Range("I7").Value = 38

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I2").Value = 24
$ws.Range("I7").Value = 101

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I3").Value = 12
$ws.Range("I7").Value = 38

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I6").Value = 50
$ws.Range("I7").Value = 141

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I2").Value = 19
$ws.Range("I7").Value = 103

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("I2").Value = 12
$ws.Range("I7").Value = 35

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("I2").Value = 17
$ws.Range("I6").Value = 8
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I2").Value = 13
$ws.Range("I7").Value = 56

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I2").Value = 25
$ws.Range("I3").Value = 33
$ws.Range("I7").Value = 124

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("I2").Value = 8
$ws.Range("I3").Value = 10
$ws.Range("I7").Value = 30

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("I2").Value = 10
$ws.Range("I7").Value = 31

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I2").Value = 7
$ws.Range("I7").Value = 41

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 22

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("I3").Value = 10
$ws.Range("I7").Value = 21

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I3").Value = 112
$ws.Range("I7").Value = 378

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("I2").Value = 17
$ws.Range("I7").Value = 45

$ws = $wb.Worksheets.Item("Sauganash,Forest Glen")
$ws.Range("I2").Value = 5
$ws.Range("I6").Value = 9
